$wb = $excel.ActiveWorkbook

# Add the new worksheet "Translations question" right after "Translations"
$ws1 = $wb.Worksheets.Item("Translations")
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Translations question"

# Header row
$ws2.Range("A1").Value = "Entity Id"
$ws2.Range("B1").Value = "Type"
$ws2.Range("C1").Value = "Index"
$ws2.Range("D1").Value = "Original"
$ws2.Range("E1").Value = "Translation"

# Data row
$ws2.Range("B2").Value = "OptionTitle"
$ws2.Range("C2").Value = 1
$ws2.Range("D2").Value = "Combobox Option"
$ws2.Range("E2").Value = "Опция Комбобокса"

# column A holds a long numeric-looking string, must be stored/kept as text
$ws2.Range("A2").NumberFormat = "@"
$ws2.Range("A2").Value = "11111111111111111111111111111111"

# Column widths
$ws2.Columns.Item(1).ColumnWidth = 43.28515625
$ws2.Columns.Item(2).ColumnWidth = 11.140625
$ws2.Columns.Item(3).ColumnWidth = 11.42578125
$ws2.Columns.Item(4).ColumnWidth = 16.85546875
$ws2.Columns.Item(5).ColumnWidth = 18.42578125

# Selection bookkeeping: new sheet becomes the active one, selection on A3
$ws2.Range("A3").Select()
